$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "TrunkROM" column header to "TRUNK ROM"
$ws.Range("C1").Value = "TRUNK ROM"

# Add new rows of ROM data (subjects 7, 8 and 9)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 24
$ws.Range("D8").Value = 34.546051025390618

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 22
$ws.Range("D9").Value = 23.148506164550781

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 38
$ws.Range("C10").Value = 35.719680665195241
